# Update scraped "want to go" counts (column F) and one ticket price (column G)
# on the two worksheets that carry the full event list: "展览" (sheet 1) and
# "全部类型" (sheet 4). "演出" and "本地生活" are untouched by this refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 2133
$ws1.Range("F9").Value  = 8647
$ws1.Range("F10").Value = 8647
$ws1.Range("F14").Value = 5786
$ws1.Range("G14").Value = 70
$ws1.Range("F16").Value = 2777
$ws1.Range("F21").Value = 622
$ws1.Range("F22").Value = 92
$ws1.Range("F23").Value = 3951
$ws1.Range("F26").Value = 64
$ws1.Range("F30").Value = 5588
$ws1.Range("F31").Value = 13
$ws1.Range("F32").Value = 73
$ws1.Range("F35").Value = 165
$ws1.Range("F36").Value = 403
$ws1.Range("F37").Value = 2764
$ws1.Range("F41").Value = 5156
$ws1.Range("F44").Value = 48
$ws1.Range("F45").Value = 3634

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 2133
$ws4.Range("F9").Value  = 8647
$ws4.Range("F10").Value = 8647
$ws4.Range("F13").Value = 5786
$ws4.Range("G13").Value = 70
$ws4.Range("F15").Value = 2777
$ws4.Range("F21").Value = 622
$ws4.Range("F22").Value = 92
$ws4.Range("F23").Value = 3951
$ws4.Range("F26").Value = 64
$ws4.Range("F30").Value = 5588
$ws4.Range("F31").Value = 13
$ws4.Range("F32").Value = 73
$ws4.Range("F34").Value = 165
$ws4.Range("F35").Value = 403
$ws4.Range("F37").Value = 2764
$ws4.Range("F42").Value = 5156
$ws4.Range("F45").Value = 48
$ws4.Range("F46").Value = 3634
